# chat/interview/questions.xlsx
# "fix: changed interview questions back to full sentences"
#
# The "chef" and "anställa/verksamhet" interview questions had been split
# into a leading transition fragment (column H) plus two or three
# lower-case sentence-continuation fragments (columns I/J/K). This restores
# them to standalone, capitalized, full-sentence questions.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4 - "bra chef" questions
$ws.Range("H4").Value = "Tack, det låter bra. När det gäller chefer.."
$ws.Range("I4").Value = "Hur tycker du att en bra chef ska vara?"
$ws.Range("J4").Value = "Finns det några speciella egenskaper du tycker att en bra chef ska ha?"
$ws.Range("K4").Value = "Är det något speciellt du skulle uppskatta hos en chef?"

# Row 5 - "beskriva dig själv / kollega" questions
$ws.Range("I5").Value = "Om du skulle beskriva dig själv som arbetskamrat, vad skulle du säga då?"
$ws.Range("J5").Value = "Om du skulle beskriva dig själv, hur tror du att du skulle vara som kollega?"

# Row 9 - "anställa dig / verksamhet" questions
$ws.Range("J9").Value = "Vad tänker du att du kan tillföra vår verksamhet?"
$ws.Range("K9").Value = "Finns det något speciellt som du tror att du kan bidra med till vår verksamhet?"
$ws.Range("I9").Value = "Jag skulle vilja veta varför du tycker att vi ska anställa just dig?"
$ws.Range("H9").Value = "Okej, jag förstår. Vi går vidare."

# Row 5 - remaining transition / closing question
$ws.Range("H5").Value = "Vad bra. Nu kommer jag ställa en fråga om dig som arbetskamrat."
$ws.Range("K5").Value = "Vad skulle du tro att andra tycker om dig som kollega?"

# Reflect the reviewer having scrolled back up and landed on I13.
$ws.Range("I13").Select()
